$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Insert a new column before A ----
# This shifts the old A ("URL TYPE"/"LOGIN"/"PRIMARY REGISTRATION URL") to B
# and the old B ("URL"/the two URLs) to C, carrying their bestFit widths along.
$ws.Columns.Item(1).Insert()

# ---- Row 1 (headers) ----
$ws.Range("A1").Value = "URL TYPE"
$ws.Range("D1").Value = "PARAMETERS"
$ws.Range("B1").Value = "URL Name"
$ws.Range("C1").Value = "URL"
$ws.Range("A1").Style = "Normal"
$ws.Range("D1").Style = "Normal"
$ws.Range("A1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true

# ---- Row 2 ----
$ws.Range("A2").Value = "POST"
$ws.Range("B2").Value = "LOGIN"
$ws.Range("C2").Value = "http://192.168.100.19/thaimaiapp/api/login/check/"

# ---- Row 3 col A (set now so shared-string order matches original authoring) ----
$ws.Range("A3").Value = "GET"

# ---- Row 2 cols D/E ----
$ws.Range("D2").Value = "picme_id,`ndob"
$ws.Range("E2").Value = " `"status`": 0`n`"message`":`"`""

# ---- Row 3 cols B/C ----
$ws.Range("B3").Value = "PRIMARY REGISTRATION URL"
$ws.Range("C3").Value = "http://192.168.100.19/thaimaiapp/api/mother/mPrimaryInfo/"

# ---- Row 4 (new row) ----
$ws.Range("A4").Value = "GET"
$ws.Range("B4").Value = "DASHBOARD"
$ws.Range("C4").Value = "http://192.168.100.19/thaimaiapp/api/mother/mDashboard"

# ---- Alignment / formatting ----
# Row 2 & 3, cols A-C: centered horizontally + vertically, taller rows
$rng223 = $ws.Range("A2:C3")
$rng223.HorizontalAlignment = -4108
$rng223.VerticalAlignment = -4108

# Row 2, cols D-E: centered + wrap text
$rngDE2 = $ws.Range("D2:E2")
$rngDE2.HorizontalAlignment = -4108
$rngDE2.VerticalAlignment = -4108
$rngDE2.WrapText = $true

# Row 4: centered horizontally only
$rng4 = $ws.Range("A4:C4")
$rng4.HorizontalAlignment = -4108

# Row heights
$ws.Rows.Item(2).RowHeight = 58.5
$ws.Rows.Item(3).RowHeight = 58.5

# ---- Hyperlink on C3 ----
$ws.Hyperlinks.Add($ws.Range("C3"), "http://192.168.100.19/thaimaiapp/api/mother/mPrimaryInfo/")

# ---- Column widths for the new columns ----
$ws.Columns.Item(1).ColumnWidth = 21.58334
$ws.Columns.Item(4).ColumnWidth = 31.75
$ws.Columns.Item(5).ColumnWidth = 17.75

# ---- Selection matches target ----
$ws.Range("C4").Select() | Out-Null
